$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds "Price" values stored as plain text in the source workbook
# (e.g. "38.617.28", "14.84"). Some of the new prices are plain decimals like
# "228.71", which Excel would otherwise silently reinterpret as a number when
# assigned through .Value. For those cells we force NumberFormat "@" (Text)
# first so the text is preserved exactly, matching the original inlineStr cell.
# Column E ("Volume(1h)") percentages such as "  +1.89%  " are never
# auto-parsed as numbers by Excel, so no special handling is required there.

$ws.Range("D2").Value = '38.617.28'
$ws.Range("E2").Value = '  +1.89%  '

$ws.Range("D3").Value = '2.092.39'
$ws.Range("E3").Value = '  +2.42%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.71'
$ws.Range("E5").Value = '  +0.01%  '

$ws.Range("E6").Value = '  +0.72%  '

$ws.Range("E7").Value = '  +0.60%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  +1.22%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0841'
$ws.Range("E10").Value = '  +2.47%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.104'
$ws.Range("E11").Value = '  +0.04%  '

$ws.Range("D12").Value = '2.398.49'
$ws.Range("E12").Value = '  +2.32%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.84'
$ws.Range("E13").Value = '  +0.50%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.31'
$ws.Range("E14").Value = '  +5.89%  '

$ws.Range("E15").Value = '  +0.26%  '

$ws.Range("E16").Value = '  +4.94%  '

$ws.Range("D17").Value = '2.093.50'
$ws.Range("E17").Value = '  +2.47%  '

$ws.Range("D18").Value = '38.513.08'
$ws.Range("E18").Value = '  +1.79%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.08'
$ws.Range("E19").Value = '  +2.76%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.98'
$ws.Range("E20").Value = '  +1.98%  '

$ws.Range("D21").Value = [string]::Concat('0.0', [char]0x2083, '0836')
$ws.Range("E21").Value = '  +1.34%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.25'
$ws.Range("E22").Value = '  +1.01%  '

$ws.Range("E23").Value = '  +0.11%  '

$ws.Range("E24").Value = '  +1.78%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.31'
$ws.Range("E25").Value = '  +1.31%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.33'
$ws.Range("E26").Value = '  +1.36%  '

$ws.Range("E27").Value = '  +0.54%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.137'
$ws.Range("E28").Value = '  +6.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.12'
$ws.Range("E29").Value = '  +1.44%  '

$ws.Range("E30").Value = '  +6.66%  '

$ws.Range("E31").Value = '  -0.46%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.34'
$ws.Range("E32").Value = '  +5.33%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.75'
$ws.Range("E33").Value = '  +4.84%  '

$ws.Range("E34").Value = '  +2.18%  '

$ws.Range("E35").Value = '  +0.43%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.53'
$ws.Range("E36").Value = '  +0.39%  '

$ws.Range("E37").Value = '  +1.68%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.56'
$ws.Range("E38").Value = '  +4.17%  '

$ws.Range("E39").Value = '  +0.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.49'
$ws.Range("E40").Value = '  +2.22%  '

$ws.Range("D41").Value = '1.545.17'
$ws.Range("E41").Value = '  +0.67%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.91'
$ws.Range("E42").Value = '  +3.86%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0221'
$ws.Range("E43").Value = '  +1.61%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.83'
$ws.Range("E44").Value = '  +0.78%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0924'
$ws.Range("E45").Value = '  +0.92%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.81'
$ws.Range("E46").Value = '  +10.67%  '

$ws.Range("E47").Value = '  +2.48%  '

$ws.Range("E48").Value = '  +1.11%  '

$ws.Range("E49").Value = '  +2.27%  '

$ws.Range("E50").Value = '  +0.28%  '

$ws.Range("D51").Value = '2.289.48'
$ws.Range("E51").Value = '  +2.57%  '
